# Update "想去人数" (number of people interested) counts for a few events
# that appear on multiple sheets ("展览" and "全部类型").

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Worksheet 1): F4 1724->1731, F5 768->769, F6 192->193
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1731
$ws1.Range("F5").Value = 769
$ws1.Range("F6").Value = 193

# Sheet "全部类型" (Worksheet 4): F4 1724->1731, F6 768->769, F7 192->193
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1731
$ws4.Range("F6").Value = 769
$ws4.Range("F7").Value = 193
